$wb = $excel.ActiveWorkbook

# --- Sheet: Q1_20_21 ---
$ws1 = $wb.Worksheets.Item("Q1_20_21")
$ws1.Range("B4").Value = "SoT"
$ws1.Range("I4").Value = "928 -678"
$ws1.Range("B5").Value = "A13"
$ws1.Range("B6").Value = "F9"

# --- Sheet: Q4_19_20 ---
$ws2 = $wb.Worksheets.Item("Q4_19_20")
$ws2.Range("B4").Value = "SoT"
$ws2.Range("B5").Value = "A11"
$ws2.Range("B6").Value = "A13"
$ws2.Range("B7").Value = "F9"

# --- Sheet: Count ---
$ws3 = $wb.Worksheets.Item("Count")
$ws3.Range("C7").Value = 833
$ws3.Range("C11").Value = 7050
